# Insert a new data row at row 253 (pushing the existing rows 253-354 down
# to 254-355) and populate it with the new weekly price observation.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows.Item(253).Insert()

$ws.Range("A253").Value = 10
$ws.Range("B253").Value = "Vega Modelo de Temuco"
$ws.Range("C253").Value = "La Araucanía"
$ws.Range("D253").Value = 44924
$ws.Range("E253").Value = 9
$ws.Range("F253").Value = "Fruta"
$ws.Range("G253").Value = 100102
$ws.Range("H253").Value = "Cítricos"
$ws.Range("I253").Value = 100102006
$ws.Range("J253").Value = "Pomelo"
$ws.Range("K253").Value = "Start Ruby"
$ws.Range("L253").Value = "Primera"
$ws.Range("M253").Value = 85
$ws.Range("N253").Value = 14000
$ws.Range("O253").Value = 14000
$ws.Range("P253").Value = 14000
$ws.Range("Q253").Value = "$/bandeja 15 kilos granel"
$ws.Range("R253").Value = "Región de O'Higgins"
$ws.Range("S253").Value = 933
$ws.Range("T253").Value = 15
